$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: fill in D19 and F19 (E19 stays blank)
$ws.Range("D19").Value = 5
$ws.Range("F19").Value = 5

# Row 23: fill in E23 and F23 (D23 already had a value)
$ws.Range("E23").Value = 5
$ws.Range("F23").Value = 5

# Row 30: fill in D30, E30 and F30
$ws.Range("D30").Value = 5
$ws.Range("E30").Value = 5
$ws.Range("F30").Value = 5

# The user scrolled the frozen window down so row 10 is the first
# visible scrollable row, then selected G30.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G30").Select()
